# lsh_coding.xlsx — add custom ordinal "order" columns to the age_groups
# sheet so that downstream output tables can be split/sorted by any of the
# three age-group schemes (official decade bands, three-way band, simple
# 0-50/51+ band) and not just age_group_simple.

$wb  = $excel.ActiveWorkbook

# --- priority_categories sheet: widen columns B and C that hold the longer
#     category labels now that age_groups feeds more output tables ---
$ws7 = $wb.Worksheets.Item("priority_categories")
$ws7.Columns.Item(2).ColumnWidth = 18.498697916666668   # -> stored width 19.33203125
$ws7.Columns.Item(3).ColumnWidth = 20.666666666666668   # -> stored width 21.5
$ws7.Range("A5").Select()

$ws9 = $wb.Worksheets.Item("age_groups")

# --- helpers: recreate the three age-group bucketing rules already used
#     in columns B (age_group_official), C (age_group_three) and
#     D (age_group_simple), then rank each bucket so it can be sorted /
#     split on numerically. ---

$officialOrder = @{
    "0-9"   = 1
    "10-19" = 2
    "20-29" = 3
    "30-39" = 4
    "40-49" = 5
    "50-59" = 6
    "60-69" = 7
    "70-79" = 8
    "80+"   = 9
}

$threeOrder = @{
    "0-50"  = 1
    "51-74" = 2
    "75+"   = 3
}

$simpleOrder = @{
    "0-50" = 1
    "51+"  = 2
}

function Get-OfficialGroup($age) {
    $lo = [Math]::Floor($age / 10) * 10
    if ($lo -ge 80) { return "80+" }
    return "$lo-$($lo + 9)"
}

function Get-ThreeGroup($age) {
    if ($age -le 50) { return "0-50" }
    elseif ($age -le 74) { return "51-74" }
    else { return "75+" }
}

function Get-SimpleGroup($age) {
    if ($age -le 50) { return "0-50" }
    else { return "51+" }
}

# --- new headers: E = age_group_order_official, F = age_group_order_three,
#     G = age_group_order_simple ---
$ws9.Cells.Item(1, 5).Value = "age_group_order_official"
$ws9.Cells.Item(1, 6).Value = "age_group_order_three"
$ws9.Cells.Item(1, 7).Value = "age_group_order_simple"

# --- fill in the ordinal rank for every age row (rows 2..122 => age 0..120) ---
for ($row = 2; $row -le 122; $row++) {
    $age = $row - 2

    $official = Get-OfficialGroup $age
    $three    = Get-ThreeGroup $age
    $simple   = Get-SimpleGroup $age

    $ws9.Cells.Item($row, 5).Value = $officialOrder[$official]
    $ws9.Cells.Item($row, 6).Value = $threeOrder[$three]
    $ws9.Cells.Item($row, 7).Value = $simpleOrder[$simple]
}

# widen the three new columns so the longer header text is visible
$ws9.Columns.Item(5).ColumnWidth = 20.830729166666668   # -> stored width 21.6640625
$ws9.Columns.Item(6).ColumnWidth = 19.498697916666668   # -> stored width 20.33203125
$ws9.Columns.Item(7).ColumnWidth = 19.666666666666668   # -> stored width 20.5

# print as a single portrait page
$ws9.PageSetup.PaperSize = 9
$ws9.PageSetup.Orientation = 1

# scroll down / select the newly added column G data that was just reviewed
# (age_groups stays the active sheet/tab, as it was before the edit)
$ws9.Range("G53:G122").Select()
